# Swap the presentation's applied colour theme from the custom "Integral"
# theme to the default "Office Theme" palette (ppt/theme/theme1.xml).
#
# COM equivalent of editing <a:clrScheme> inside ppt/theme/theme1.xml:
# PowerPoint exposes the twelve theme colours (dk1, lt1, dk2, lt2,
# accent1-6, hlink, folHlink - in that fixed order) through
# Theme.ThemeColorScheme.Item(1..12).RGB, reachable from the slide
# master that owns theme1.xml.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$theme = $master.Theme
$colors = $theme.ThemeColorScheme

# New values = the stock "Office Theme" palette (RGB hex -> COM long,
# long = R + G*256 + B*65536).
$officeTheme = @(
    0,         # 1  dk1      000000
    16777215,  # 2  lt1      FFFFFF
    6968388,   # 3  dk2      44546A
    15132391,  # 4  lt2      E7E6E6
    13998939,  # 5  accent1  5B9BD5
    3243501,   # 6  accent2  ED7D31
    10855845,  # 7  accent3  A5A5A5
    49407,     # 8  accent4  FFC000
    12874308,  # 9  accent5  4472C4
    4697456,   # 10 accent6  70AD47
    12673797,  # 11 hlink    0563C1
    7491477    # 12 folHlink 954F72
)

for ($i = 1; $i -le $colors.Count; $i++) {
    $colors.Item($i).RGB = $officeTheme[$i - 1]
}
